# Update gh-pages to output generated at 456a3b4
# The scraper dropped the oldest (already-elapsed) event from each listing
# sheet and refreshed the "want to go" counters (column F) for the events
# that remain. This affects the "展览" (sheet 1) and "全部类型" (sheet 4)
# tabs; "演出" and "本地生活" are untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# The first data row (row 2, "南宁·AB动漫游戏嘉年华") has already
# happened and drops off the list; deleting it shifts every other
# row's B:I content up by one.
$ws1.Range("A2:I2").EntireRow.Delete()

# Column A is a static running index (0,1,2,...), not a formula, so the
# delete leaves it one-off; put the 1..N sequence back.
$dims1 = $ws1.UsedRange.Rows.Count
for ($r = 2; $r -le $dims1; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}

# Refresh the "想去人数" (want-to-go) counters that ticked up between
# scrapes.
$ws1.Cells.Item(2, 6).Value = 13101
$ws1.Cells.Item(3, 6).Value = 316
$ws1.Cells.Item(4, 6).Value = 639
$ws1.Cells.Item(5, 6).Value = 209
$ws1.Cells.Item(6, 6).Value = 413
$ws1.Cells.Item(7, 6).Value = 1243
$ws1.Cells.Item(8, 6).Value = 119

# ---------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("A2:I2").EntireRow.Delete()

$dims4 = $ws4.UsedRange.Rows.Count
for ($r = 2; $r -le $dims4; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

$ws4.Cells.Item(3, 6).Value = 13101
$ws4.Cells.Item(4, 6).Value = 316
$ws4.Cells.Item(5, 6).Value = 639
$ws4.Cells.Item(6, 6).Value = 209
$ws4.Cells.Item(9, 6).Value = 413
$ws4.Cells.Item(10, 6).Value = 1243
$ws4.Cells.Item(12, 6).Value = 119
